$wb = $excel.ActiveWorkbook

# ALC row 15: Morning Glass of Ether
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 683.61
$ws.Range("I15").Value = 683.61
$ws.Range("K15").Value = 2050.83
$ws.Range("M15").Value = -1881.83

# ALC row 17: One for the Road
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2865114.8
$ws.Range("J17").Value = 2917206
$ws.Range("L17").Value = 8751618
$ws.Range("N17").Value = -8751954

# ALC row 53: No Accounting for Waste
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 35786620
$ws.Range("I53").Value = 50101084
$ws.Range("J53").Value = 469.75
$ws.Range("K53").Value = 50101084
$ws.Range("L53").Value = 469.75
$ws.Range("M53").Value = -50100447
$ws.Range("N53").Value = -1743.75

# ALC row 98: The Dotted Line
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 947.3333
$ws.Range("I98").Value = 942.5454999999999
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 942.5454999999999
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 555.4545000000001
$ws.Range("N98").Value = -3996

# ALC row 122: Wishful Inking
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 947.3333
$ws.Range("I122").Value = 942.5454999999999
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2827.6365
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -377.6364999999996
$ws.Range("N122").Value = -7900

# ALC row 129: Practical Command
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 914.087
$ws.Range("J129").Value = 1003.7455
$ws.Range("L129").Value = 3011.2365
$ws.Range("N129").Value = -13011.2365

# ARM row 2: Ain't Got No Ingots
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2664.3635
$ws.Range("I2").Value = 1977.8572
$ws.Range("J2").Value = 3865.75
$ws.Range("K2").Value = 1977.8572
$ws.Range("L2").Value = 3865.75
$ws.Range("M2").Value = -1864.8572
$ws.Range("N2").Value = -4091.75

# ARM row 32: Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4681.883
$ws.Range("I32").Value = 3425.3164
$ws.Range("J32").Value = 11299.8
$ws.Range("K32").Value = 3425.3164
$ws.Range("L32").Value = 11299.8
$ws.Range("M32").Value = -3138.3164
$ws.Range("N32").Value = -11873.8

# ARM row 116: No Scope
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2664.3635
$ws.Range("I116").Value = 1977.8572
$ws.Range("J116").Value = 3865.75
$ws.Range("K116").Value = 1977.8572
$ws.Range("L116").Value = 3865.75
$ws.Range("M116").Value = 316.1428000000001
$ws.Range("N116").Value = -8453.75

# ARM row 122: Haste for High Durium
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2335822.8
$ws.Range("I122").Value = 3210668.5
$ws.Range("K122").Value = 9632005.5
$ws.Range("M122").Value = -9629555.5

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2616.1304
$ws.Range("I132").Value = 1712.8379
$ws.Range("K132").Value = 5138.5137
$ws.Range("M132").Value = -2608.5137

# BSM row 3: Hells Bells
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2664.3635
$ws.Range("I3").Value = 1977.8572
$ws.Range("J3").Value = 3865.75
$ws.Range("K3").Value = 1977.8572
$ws.Range("L3").Value = 3865.75
$ws.Range("M3").Value = -1863.8572
$ws.Range("N3").Value = -4093.75

# BSM row 20: Smelt and Dealt
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43118
$ws.Range("I20").Value = 1602.6666
$ws.Range("J20").Value = 84633.336
$ws.Range("K20").Value = 1602.6666
$ws.Range("L20").Value = 84633.336
$ws.Range("M20").Value = -1355.6666
$ws.Range("N20").Value = -85127.336

# BSM row 99: Meddle in Metal
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 90910260
$ws.Range("I99").Value = 125001110
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 125001110
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = -124999612
$ws.Range("N99").Value = -4296

# BSM row 105: Ingot to Wing It
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 22728832
$ws.Range("I105").Value = 27779012
$ws.Range("J105").Value = 3024.75
$ws.Range("K105").Value = 27779012
$ws.Range("L105").Value = 3024.75
$ws.Range("M105").Value = -27777265
$ws.Range("N105").Value = -6518.75

# BSM row 107: The Gold Experience
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 71429496
$ws.Range("I107").Value = 111112024
$ws.Range("J107").Value = 949.2
$ws.Range("K107").Value = 111112024
$ws.Range("L107").Value = 949.2
$ws.Range("M107").Value = -111110104
$ws.Range("N107").Value = -4789.2

# BSM row 134: Ruthenium Supremium
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2805.7808
$ws.Range("I134").Value = 2882.3667
$ws.Range("K134").Value = 8647.1001
$ws.Range("M134").Value = -6112.1001

# CRP row 31: Wall Not Found
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5370817.5
$ws.Range("I31").Value = 1465.1666
$ws.Range("J31").Value = 12274270
$ws.Range("K31").Value = 1465.1666
$ws.Range("L31").Value = 12274270
$ws.Range("M31").Value = -1170.1666
$ws.Range("N31").Value = -12274860

# CRP row 34: Armoires of the Rich and Famous
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5370817.5
$ws.Range("I34").Value = 1465.1666
$ws.Range("J34").Value = 12274270
$ws.Range("K34").Value = 1465.1666
$ws.Range("L34").Value = 12274270
$ws.Range("M34").Value = -1263.1666
$ws.Range("N34").Value = -12274674

# CRP row 58: You Do the Heavy Lifting
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2826043.8
$ws.Range("I58").Value = 4505220
$ws.Range("J58").Value = 1974.5454
$ws.Range("K58").Value = 4505220
$ws.Range("L58").Value = 1974.5454
$ws.Range("M58").Value = -4505017
$ws.Range("N58").Value = -2380.5454

# CRP row 87: Anatomy of a Drill Bit
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = ""
$ws.Range("N87").Value = 0

# CRP row 90: Pulling Them to the Grind (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = ""
$ws.Range("N90").Value = 0

# CRP row 122: Timber of Tenkonto
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3379.2222
$ws.Range("I122").Value = 2568.8333
$ws.Range("K122").Value = 7706.499899999999
$ws.Range("M122").Value = -5256.499899999999

# CRP row 134: Wood You Be Quiet
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9011454
$ws.Range("I134").Value = 15876389
$ws.Range("J134").Value = 1228.125
$ws.Range("K134").Value = 47629167
$ws.Range("L134").Value = 3684.375
$ws.Range("M134").Value = -47626632
$ws.Range("N134").Value = -8754.375

# CRP row 136: Turali Quality
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2826043.8
$ws.Range("I136").Value = 4505220
$ws.Range("J136").Value = 1974.5454
$ws.Range("K136").Value = 13515660
$ws.Range("L136").Value = 5923.6362
$ws.Range("M136").Value = -13513110
$ws.Range("N136").Value = -11023.6362

# CUL row 107: Slippery Service
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4530.68
$ws.Range("I107").Value = 365.5
$ws.Range("J107").Value = 5846
$ws.Range("K107").Value = 1096.5
$ws.Range("L107").Value = 17538
$ws.Range("M107").Value = 823.5
$ws.Range("N107").Value = -21378

# GSM row 122: Awarding Academic Excellence
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 65579324
$ws.Range("I122").Value = 133103280
$ws.Range("J122").Value = 5558033
$ws.Range("K122").Value = 399309840
$ws.Range("L122").Value = 16674099
$ws.Range("M122").Value = -399307390
$ws.Range("N122").Value = -16678999

# GSM row 126: Gold Rush Order
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4496.2793
$ws.Range("I126").Value = 6581.476
$ws.Range("J126").Value = 2505.8635
$ws.Range("K126").Value = 19744.428
$ws.Range("L126").Value = 7517.5905
$ws.Range("M126").Value = -17274.428
$ws.Range("N126").Value = -12457.5905

# LTW row 132: Tenets of Tanning
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 14820005
$ws.Range("I132").Value = 18396248
$ws.Range("K132").Value = 55188744
$ws.Range("M132").Value = -55186214

# WVR row 122: Heavy Armoire
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1409.5264
$ws.Range("I122").Value = 1305.0625
$ws.Range("K122").Value = 3915.1875
$ws.Range("M122").Value = -1465.1875

# WVR row 136: Weaving the Envelope
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 12823221
$ws.Range("I136").Value = 3014.6572
$ws.Range("J136").Value = 39217764
$ws.Range("K136").Value = 9043.971600000001
$ws.Range("L136").Value = 117653292
$ws.Range("M136").Value = -6493.971600000001
$ws.Range("N136").Value = -117658392
